# Added ACMG Tier3 data
# - Rename "Strategy A2 (FS)" header to "Strategy A (FS)"
# - Add new "Strategy B (ACOG+ACMG)" column (G) with per-disease flag data
# - Resize the data columns (C:F) to fit the new layout
# - Move the active selection to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the existing "Strategy A2 (FS)" header text -----------------------
$ws.Range("D1").Value = "Strategy A (FS)"

# --- New header for the added strategy column -------------------------------
$ws.Range("G1").Value = "Strategy B (ACOG+ACMG)"

# --- New column G data (ACMG Tier3 flags), rows 2-22 ------------------------
$g = @(1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $g.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $g[$i]
}

# --- Column widths for C, D, E, F, to fit the new layout --------------------
$ws.Columns.Item(3).ColumnWidth = 15.57 - 5/6
$ws.Columns.Item(4).ColumnWidth = 13.93 - 5/6
$ws.Columns.Item(5).ColumnWidth = 16.79 - 5/6
$ws.Columns.Item(6).ColumnWidth = 17.13 - 5/6

# --- Move the active selection --------------------------------------------
$ws.Range("D2").Select()
